# Development Strategy.docx update
#  - swap the two inline pictures' display names (image1.png <-> image2.png)
#    (the embedded picture data/relationships themselves are unchanged)
#  - bump the "Last modified/reviewed" dates 06/30/2022 -> 06/27/2023
#  - fix "JsonWeb token" -> "JSON Web token"
#  - expand the JWT service description
#  - fix the "Deply: " typo -> "Deployment: "

$d = $word.ActiveDocument

# --- text fixes (table cells + bullet) -------------------------------------

$d.Content.Find.Execute(
    "06/30/2022", $true, $false, $false, $false, $false, $true, 1, $false,
    "06/27/2023", 2) | Out-Null

$d.Content.Find.Execute(
    "JsonWeb token", $true, $false, $false, $false, $false, $true, 1, $false,
    "JSON Web token", 2) | Out-Null

$d.Content.Find.Execute(
    "Service that provides secure access, as well as account management, roles and permissions.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Service that provides secure access, for the creation of access tokens that allow the propagation of identity and privileges",
    2) | Out-Null

$d.Content.Find.Execute(
    "Deply: ", $true, $false, $false, $false, $false, $true, 1, $false,
    "Deployment: ", 2) | Out-Null

# --- swap the two inline pictures' displayed names --------------------------
# InlineShapes have no direct "Name" setter on the Word object model (unlike
# floating Shapes), so round-trip the canonical OOXML and swap the
# wp:docPr / pic:cNvPr "name" attributes there. A placeholder token avoids a
# naive replace-then-replace undoing itself.

$xml = $d.WordOpenXML
$placeholder = 'image1.png'
$other = 'image2.png'

$xml = $xml -replace [regex]::Escape('name="image1.png"'), 'name="__TMP_IMAGE_NAME_SWAP__"'
$xml = $xml -replace [regex]::Escape('name="image2.png"'), 'name="image1.png"'
$xml = $xml -replace [regex]::Escape('name="__TMP_IMAGE_NAME_SWAP__"'), 'name="image2.png"'

$d.WordOpenXML = $xml
